$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Swap the contents of row 19 and row 20 (all columns changed between
# the two species records), including moving the "Publik kommentar"
# (AC) value from row 19 to row 20.

$row19_A = $ws.Range("A19").Value2
$row19_B = $ws.Range("B19").Value2
$row19_E = $ws.Range("E19").Value2
$row19_F = $ws.Range("F19").Value2
$row19_G = $ws.Range("G19").Value2
$row19_H = $ws.Range("H19").Value2
$row19_Q = $ws.Range("Q19").Value2
$row19_R = $ws.Range("R19").Value2
$row19_AC = $ws.Range("AC19").Value2

$row20_A = $ws.Range("A20").Value2
$row20_B = $ws.Range("B20").Value2
$row20_E = $ws.Range("E20").Value2
$row20_F = $ws.Range("F20").Value2
$row20_G = $ws.Range("G20").Value2
$row20_H = $ws.Range("H20").Value2
$row20_Q = $ws.Range("Q20").Value2
$row20_R = $ws.Range("R20").Value2

$ws.Range("A19").Value2 = $row20_A
$ws.Range("B19").Value2 = $row20_B
$ws.Range("E19").Value2 = $row20_E
$ws.Range("F19").Value2 = $row20_F
$ws.Range("G19").Value2 = $row20_G
$ws.Range("H19").Value2 = $row20_H
$ws.Range("Q19").Value2 = $row20_Q
$ws.Range("R19").Value2 = $row20_R
$ws.Range("AC19").Value2 = $null

$ws.Range("A20").Value2 = $row19_A
$ws.Range("B20").Value2 = $row19_B
$ws.Range("E20").Value2 = $row19_E
$ws.Range("F20").Value2 = $row19_F
$ws.Range("G20").Value2 = $row19_G
$ws.Range("H20").Value2 = $row19_H
$ws.Range("Q20").Value2 = $row19_Q
$ws.Range("R20").Value2 = $row19_R
$ws.Range("AC20").Value2 = $row19_AC
